$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Palermo Soho'
$ws.Cells.Item(2, 2).Value = 'USD 135.000'
$ws.Cells.Item(2, 3).Value = '2 Ambientes - Apto Profesional - Dueño Directo'
$ws.Cells.Item(2, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-palermo-dueno-directo.html'

$ws.Cells.Item(3, 1).Value = 'Recoleta'
$ws.Cells.Item(3, 2).Value = 'USD 110.000'
$ws.Cells.Item(3, 3).Value = 'Estilo Francés - Sin Comisión Inmobiliaria'
$ws.Cells.Item(3, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-recoleta-dueno-directo.html'

$ws.Cells.Item(4, 1).Value = 'Belgrano R'
$ws.Cells.Item(4, 2).Value = 'USD 145.000'
$ws.Cells.Item(4, 3).Value = '3 Ambientes con Cochera - Dueño Vende'
$ws.Cells.Item(4, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-belgrano-dueno-directo.html'

$ws.Cells.Item(5, 1).Value = 'Caballito Centro'
$ws.Cells.Item(5, 2).Value = 'USD 88.000'
$ws.Cells.Item(5, 3).Value = 'Oportunidad Retasado - Dueño Directo'
$ws.Cells.Item(5, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-caballito-dueno-directo.html'

$ws.Cells.Item(6, 1).Value = 'Villa Urquiza'
$ws.Cells.Item(6, 2).Value = 'USD 105.000'
$ws.Cells.Item(6, 3).Value = 'Monoambiente Divisible - Estreno - S/Comisión'
$ws.Cells.Item(6, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-villa-urquiza-dueno-directo.html'

$ws.Cells.Item(7, 1).Value = 'Almagro'
$ws.Cells.Item(7, 2).Value = 'USD 72.000'
$ws.Cells.Item(7, 3).Value = 'Ideal Inversión Rentabilidad 5% anual'
$ws.Cells.Item(7, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-almagro-dueno-directo.html'

$ws.Cells.Item(8, 1).Value = 'Nuñez'
$ws.Cells.Item(8, 2).Value = 'USD 128.000'
$ws.Cells.Item(8, 3).Value = 'Cerca del Subte D - Dueño Directo Impecable'
$ws.Cells.Item(8, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-nunez-dueno-directo.html'

$ws.Cells.Item(9, 1).Value = 'Flores'
$ws.Cells.Item(9, 2).Value = 'USD 65.000'
$ws.Cells.Item(9, 3).Value = '2 Ambientes Luminoso - Oportunidad Efectivo'
$ws.Cells.Item(9, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-flores-dueno-directo.html'

$ws.Cells.Item(10, 1).Value = 'Villa Crespo'
$ws.Cells.Item(10, 2).Value = 'USD 92.000'
$ws.Cells.Item(10, 3).Value = 'Zona Outlets - Excelente Ubicación - S/Comisión'
$ws.Cells.Item(10, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-villa-crespo-dueno-directo.html'

$ws.Cells.Item(11, 1).Value = 'San Telmo'
$ws.Cells.Item(11, 2).Value = 'USD 78.000'
$ws.Cells.Item(11, 3).Value = 'Casco Histórico - Ideal AirBnb - Dueño Directo'
$ws.Cells.Item(11, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-san-telmo-dueno-directo.html'

$ws.Cells.Item(12, 1).Value = 'Colegiales'
$ws.Cells.Item(12, 2).Value = 'USD 115.000'
$ws.Cells.Item(12, 3).Value = '3 Ambientes Amplio - Dueño Vende Urgente'
$ws.Cells.Item(12, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-colegiales-dueno-directo.html'

$ws.Cells.Item(13, 1).Value = 'Barracas'
$ws.Cells.Item(13, 2).Value = 'USD 82.000'
$ws.Cells.Item(13, 3).Value = 'Edificio Moderno - Seguridad - Sin Comisión'
$ws.Cells.Item(13, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-barracas-dueno-directo.html'

$ws.Cells.Item(14, 1).Value = 'Chacarita'
$ws.Cells.Item(14, 2).Value = 'USD 98.000'
$ws.Cells.Item(14, 3).Value = 'Punto Estratégico - 2 Ambientes Estreno'
$ws.Cells.Item(14, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-chacarita-dueno-directo.html'

$ws.Cells.Item(15, 1).Value = 'Villa Devoto'
$ws.Cells.Item(15, 2).Value = 'USD 140.000'
$ws.Cells.Item(15, 3).Value = 'Residencial - 3 Ambientes con Balcón Terraza'
$ws.Cells.Item(15, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-villa-devoto-dueno-directo.html'

$ws.Cells.Item(16, 1).Value = 'Saavedra'
$ws.Cells.Item(16, 2).Value = 'USD 108.000'
$ws.Cells.Item(16, 3).Value = 'Frente al Parque - Dueño Directo - Muy Luminoso'
$ws.Cells.Item(16, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-saavedra-dueno-directo.html'

$ws.Cells.Item(17, 1).Value = 'Balvanera'
$ws.Cells.Item(17, 2).Value = 'USD 58.000'
$ws.Cells.Item(17, 3).Value = 'Económico - Cerca de Facultades - Ideal Estudiantes'
$ws.Cells.Item(17, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-balvanera-dueno-directo.html'

$ws.Cells.Item(18, 1).Value = 'Boedo'
$ws.Cells.Item(18, 2).Value = 'USD 74.000'
$ws.Cells.Item(18, 3).Value = 'Tradicional - 2 Ambientes - Dueño Directo'
$ws.Cells.Item(18, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-boedo-dueno-directo.html'

$ws.Cells.Item(19, 1).Value = 'Coghlan'
$ws.Cells.Item(19, 2).Value = 'USD 122.000'
$ws.Cells.Item(19, 3).Value = 'Zona Tranquila - Edificio de Categoría'
$ws.Cells.Item(19, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-coghlan-dueno-directo.html'

$ws.Cells.Item(20, 1).Value = 'Puerto Madero'
$ws.Cells.Item(20, 2).Value = 'USD 350.000'
$ws.Cells.Item(20, 3).Value = 'Lujo - Vista al Río - Dueño Vende Directo'
$ws.Cells.Item(20, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-puerto-madero-dueno-directo.html'

$ws.Cells.Item(21, 1).Value = 'Villa Luro'
$ws.Cells.Item(21, 2).Value = 'USD 87.000'
$ws.Cells.Item(21, 3).Value = 'Impecable - Sin Expensas - Dueño Directo'
$ws.Cells.Item(21, 4).Value = 'https://www.zonaprop.com.ar/departamentos-alquiler-villa-luro-dueno-directo.html'
